# Bug 2781 V 1.2.2 updates - Added CoreData stuff
#
# The V1.2.0 test-plan worksheet is duplicated (so the original V1.2.0
# content is preserved in a new tab at the end of the workbook) and the
# original tab is repurposed/updated in place to become the new
# "V1.2.2 7_10_2015" test plan (several "x" check marks are cleared and
# two path cells are generalized).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # "V1.2.0 1_26_2015 (2)" -> becomes "V1.2.2 7_10_2015"
$ws2 = $wb.Worksheets.Item(2)   # "V1.2.1 6_9_2015" (unchanged content)

# Duplicate the original V1.2.0 sheet, placing the copy at the end of the
# workbook so the original content is preserved as its own tab.
$ws1.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3 = $wb.Worksheets.Item($wb.Worksheets.Count)

# Rename sheets: the in-place sheet becomes the new V1.2.2 plan, the
# freshly duplicated sheet at the end keeps the old V1.2.0 name.
$ws1.Name = "V1.2.2 7_10_2015"
$ws3.Name = "V1.2.0 1_26_2015 (2)"

# Clear the "x" marks that are no longer applicable for the new release
# on the V1.2.2 sheet (ws1).
$ws1.Range("D5").ClearContents()
$ws1.Range("E5").ClearContents()
$ws1.Range("F5").ClearContents()
$ws1.Range("G5").ClearContents()

$ws1.Range("C6").ClearContents()
$ws1.Range("G6").ClearContents()

$ws1.Range("C7").ClearContents()

$ws1.Range("C8").ClearContents()
$ws1.Range("G8").ClearContents()

$ws1.Range("C9").ClearContents()
$ws1.Range("G9").ClearContents()

$ws1.Range("C10").ClearContents()
$ws1.Range("G10").ClearContents()

$ws1.Range("C11").ClearContents()
$ws1.Range("G11").ClearContents()

$ws1.Range("C12").ClearContents()
$ws1.Range("G12").ClearContents()

$ws1.Range("C13").ClearContents()
$ws1.Range("G13").ClearContents()

$ws1.Range("C14:G14").ClearContents()
$ws1.Range("C15:G15").ClearContents()
$ws1.Range("C16:G16").ClearContents()
$ws1.Range("C17:G17").ClearContents()

$ws1.Range("C19").ClearContents()
$ws1.Range("C20").ClearContents()
$ws1.Range("C21").ClearContents()

# Generalize the release-path comments so they no longer hard-code the
# V1.2.0 release folder name.
$ws1.Range("I19").Value = "/Users/scoleman/dev/fips/fcids/release/Vxxx.../fips-pi.xcodeproj"
$ws1.Range("I20").Value = "/Users/scoleman/dev/fips/fcids/release/Vxxx.../testcordova.xcodeproj"

# Update sheet view / selection state: V1.2.2 (ws1) becomes the active
# tab with I25 selected, V1.2.1 (ws2) keeps its own F21 selection but is
# no longer the active tab.
$ws2.Range("F21").Select()
$ws1.Activate()
$ws1.Range("I25").Select()
